{"js": "// Update the worksheet date and all the division problems to the new\n// values, per the commit's regenerated output.\nconst replacements = [\n  [\"2025-01-25 Saturday\", \"2025-01-26 Sunday\"],\n  [\"40\u00f75=\", \"30\u00f76=\"],\n  [\"43\u00f77=\", \"50\u00f77=\"],\n  [\"83\u00f78=\", \"40\u00f72=\"],\n  [\"20\u00f74=\", \"17\u00f74=\"],\n  [\"88\u00f78=\", \"15\u00f77=\"],\n  [\"33\u00f76=\", \"31\u00f74=\"],\n  [\"83\u00f76=\", \"85\u00f72=\"],\n  [\"11\u00f76=\", \"50\u00f72=\"],\n  [\"61\u00f76=\", \"20\u00f73=\"],\n  [\"56\u00f78=\", \"87\u00f77=\"],\n  [\"83\u00f79=\", \"70\u00f74=\"],\n  [\"34\u00f75=\", \"98\u00f76=\"],\n  [\"48\u00f72=\", \"63\u00f74=\"],\n  [\"29\u00f79=\", \"65\u00f75=\"],\n  [\"15\u00f75=\", \"94\u00f72=\"],\n  [\"91\u00f73=\", \"10\u00f74=\"],\n  [\"29\u00f73=\", \"31\u00f73=\"],\n  [\"94\u00f79=\", \"37\u00f78=\"],\n  [\"56\u00f76=\", \"96\u00f76=\"],\n  [\"48\u00f79=\", \"32\u00f77=\"],\n  [\"95\u00f74=\", \"47\u00f78=\"],\n  [\"57\u00f78=\", \"12\u00f78=\"],\n  [\"73\u00f76=\", \"83\u00f77=\"],\n  [\"92\u00f72=\", \"90\u00f72=\"],\n  [\"25\u00f75=\", \"18\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all the division problems to the new\n# values, per the commit's regenerated output.\n$replacements = @(\n    @(\"2025-01-25 Saturday\", \"2025-01-26 Sunday\"),\n    @(\"40\u00f75=\", \"30\u00f76=\"),\n    @(\"43\u00f77=\", \"50\u00f77=\"),\n    @(\"83\u00f78=\", \"40\u00f72=\"),\n    @(\"20\u00f74=\", \"17\u00f74=\"),\n    @(\"88\u00f78=\", \"15\u00f77=\"),\n    @(\"33\u00f76=\", \"31\u00f74=\"),\n    @(\"83\u00f76=\", \"85\u00f72=\"),\n    @(\"11\u00f76=\", \"50\u00f72=\"),\n    @(\"61\u00f76=\", \"20\u00f73=\"),\n    @(\"56\u00f78=\", \"87\u00f77=\"),\n    @(\"83\u00f79=\", \"70\u00f74=\"),\n    @(\"34\u00f75=\", \"98\u00f76=\"),\n    @(\"48\u00f72=\", \"63\u00f74=\"),\n    @(\"29\u00f79=\", \"65\u00f75=\"),\n    @(\"15\u00f75=\", \"94\u00f72=\"),\n    @(\"91\u00f73=\", \"10\u00f74=\"),\n    @(\"29\u00f73=\", \"31\u00f73=\"),\n    @(\"94\u00f79=\", \"37\u00f78=\"),\n    @(\"56\u00f76=\", \"96\u00f76=\"),\n    @(\"48\u00f79=\", \"32\u00f77=\"),\n    @(\"95\u00f74=\", \"47\u00f78=\"),\n    @(\"57\u00f78=\", \"12\u00f78=\"),\n    @(\"73\u00f76=\", \"83\u00f77=\"),\n    @(\"92\u00f72=\", \"90\u00f72=\"),\n    @(\"25\u00f75=\", \"18\u00f76=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
